# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Updates the "K" column (column G) values on Sheet1 for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 0
    3  = 1
    5  = 0
    6  = 0
    7  = 1
    8  = 2
    9  = 0
    10 = 1
    11 = 1
    12 = 0
    13 = 1
    14 = 0
    15 = 0
    16 = 0
    17 = 1
    18 = 1
    19 = 1
    21 = 2
    22 = 2
    23 = 1
    24 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
